$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing review text entries (A3, A4)
$ws.Range("A3").Value = "this is not good product"
$ws.Range("A4").Value = "nice product but costly "

# Append new review rows (A5:A7)
$ws.Range("A5").Value = "product is far away from my budget but it is very helpful "
$ws.Range("A6").Value = "your overall service is very good "
$ws.Range("A7").Value = "I am very satisfied with this product "

# Set column A width to match target layout (21.5703125 stored units).
# The ColumnWidth property here is quantized to 1/6-character-unit steps
# before Excel's internal +5/6 padding is stored in the XML, so the nearest
# reachable stored width is 21.5; 20.67 sits safely inside that rounding
# bucket.
$ws.Columns.Item(1).ColumnWidth = 20.67

# Match the new active selection cell from the diff (B8)
$ws.Range("B8").Select()
